$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '89.557.74'
Set-TextValue $ws.Range('E2') '  +3.01%  '
Set-TextValue $ws.Range('D3') '3.183.56'
Set-TextValue $ws.Range('E3') '  +1.36%  '
Set-TextValue $ws.Range('E4') '  -0.03%  '
Set-TextValue $ws.Range('D5') '216.22'
Set-TextValue $ws.Range('E5') '  +6.13%  '
Set-TextValue $ws.Range('D6') '622.95'
Set-TextValue $ws.Range('E6') '  +2.77%  '
Set-TextValue $ws.Range('D7') '0.388'
Set-TextValue $ws.Range('E7') '  +5.25%  '
Set-TextValue $ws.Range('D8') '0.688'
Set-TextValue $ws.Range('E8') '  +4.89%  '
Set-TextValue $ws.Range('E9') '  +0.08%  '
Set-TextValue $ws.Range('D10') '3.174.42'
Set-TextValue $ws.Range('E10') '  +1.46%  '
Set-TextValue $ws.Range('D11') '0.568'
Set-TextValue $ws.Range('E11') '  +7.48%  '
Set-TextValue $ws.Range('E12') '  +1.94%  '
Set-TextValue $ws.Range('D13') '0.0000256'
Set-TextValue $ws.Range('E13') '  +6.40%  '
Set-TextValue $ws.Range('D14') '5.36'
Set-TextValue $ws.Range('E14') '  +3.12%  '
Set-TextValue $ws.Range('D15') '33.27'
Set-TextValue $ws.Range('E15') '  +4.52%  '
Set-TextValue $ws.Range('D16') '3.767.38'
Set-TextValue $ws.Range('E16') '  +1.31%  '
Set-TextValue $ws.Range('D17') '89.461.76'
Set-TextValue $ws.Range('E17') '  +3.24%  '
Set-TextValue $ws.Range('D18') '3.208.41'
Set-TextValue $ws.Range('E18') '  +0.94%  '
Set-TextValue $ws.Range('D19') '3.44'
Set-TextValue $ws.Range('E19') '  +15.84%  '
Set-TextValue $ws.Range('D20') '0.0000224'
Set-TextValue $ws.Range('E20') '  +73.93%  '
Set-TextValue $ws.Range('D21') '13.43'
Set-TextValue $ws.Range('E21') '  +1.25%  '
Set-TextValue $ws.Range('D22') '433.96'
Set-TextValue $ws.Range('E22') '  +5.83%  '
Set-TextValue $ws.Range('D23') '8.60'
Set-TextValue $ws.Range('E23') '  +2.26%  '
Set-TextValue $ws.Range('D24') '5.06'
Set-TextValue $ws.Range('E24') '  +0.09%  '
Set-TextValue $ws.Range('D25') '5.28'
Set-TextValue $ws.Range('E25') '  +3.53%  '
Set-TextValue $ws.Range('D26') '11.91'
Set-TextValue $ws.Range('E26') '  +1.90%  '
Set-TextValue $ws.Range('D27') '81.67'
Set-TextValue $ws.Range('E27') '  +11.99%  '
Set-TextValue $ws.Range('D28') '3.376.57'
Set-TextValue $ws.Range('E28') '  +2.36%  '
Set-TextValue $ws.Range('D29') '0.999'
Set-TextValue $ws.Range('E29') '  +0.05%  '
Set-TextValue $ws.Range('D30') '0.158'
Set-TextValue $ws.Range('E30') '  -1.24%  '
Set-TextValue $ws.Range('D31') '1.00'
Set-TextValue $ws.Range('E31') '  -0.08%  '
Set-TextValue $ws.Range('D32') '4.04'
Set-TextValue $ws.Range('E32') '  +35.84%  '
Set-TextValue $ws.Range('B33') 'Bittensor'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D33') '542.98'
Set-TextValue $ws.Range('E33') '  +1.26%  '
Set-TextValue $ws.Range('B34') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D34') '8.43'
Set-TextValue $ws.Range('E34') '  +2.67%  '
Set-TextValue $ws.Range('D35') '7.00'
Set-TextValue $ws.Range('E35') '  +7.30%  '
Set-TextValue $ws.Range('D36') '1.90'
Set-TextValue $ws.Range('E36') '  +3.50%  '
Set-TextValue $ws.Range('D37') '1.31'
Set-TextValue $ws.Range('E37') '  +1.48%  '
Set-TextValue $ws.Range('D38') '22.27'
Set-TextValue $ws.Range('E38') '  +3.35%  '
Set-TextValue $ws.Range('D39') '22.38'
Set-TextValue $ws.Range('E39') '  +2.80%  '
Set-TextValue $ws.Range('D40') '0.127'
Set-TextValue $ws.Range('E40') '  -3.47%  '
Set-TextValue $ws.Range('D41') '0.997'
Set-TextValue $ws.Range('E41') '  +0.02%  '
Set-TextValue $ws.Range('E42') '  -0.06%  '
Set-TextValue $ws.Range('D43') '1.92'
Set-TextValue $ws.Range('E43') '  +1.57%  '
Set-TextValue $ws.Range('D44') '0.371'
Set-TextValue $ws.Range('E44') '  +1.03%  '
Set-TextValue $ws.Range('D45') '146.77'
Set-TextValue $ws.Range('E45') '  -0.83%  '
Set-TextValue $ws.Range('D46') '172.82'
Set-TextValue $ws.Range('E46') '  +1.31%  '
Set-TextValue $ws.Range('D47') '43.65'
Set-TextValue $ws.Range('E47') '  +2.04%  '
Set-TextValue $ws.Range('D48') '0.755'
Set-TextValue $ws.Range('E48') '  +10.74%  '
Set-TextValue $ws.Range('D49') '0.124'
Set-TextValue $ws.Range('E49') '  -2.02%  '
Set-TextValue $ws.Range('D50') '1.24'
Set-TextValue $ws.Range('E50') '  +0.37%  '
Set-TextValue $ws.Range('D51') '0.617'
Set-TextValue $ws.Range('E51') '  +6.14%  '
